$d = $word.ActiveDocument

function FR([string]$find, [string]$replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "MISSING: $find"
    }
}

# 1. Opening paragraph rewrite
FR "After joining Spotify last year, I was curious to know how I listened to music, so I set out to perform an exploratory analysis to answer some questions to help me understand my listing habit. " "I joined Spotify last year, I was curious to know more about my streaming habit, so I set out to perform an exploratory analysis to answer some questions to help me understand how I listen to music"

# 2. "My favorite song" -> "Top song"
FR "My favorite song" "Top song"

# 3. "It is not surprising" -> "It is not too surprising"
FR "It is not surprising to me because I love music." "It is not too surprising to me because I love music."

# 4. "surfing through Twitter." -> "surfing the web."
FR "surfing through Twitter." "surfing the web."

# 5. Business Wars / podcast(s)
FR " Naval and Business  Wars are some of the podcast I listen to; although" " Naval and Business Wars are some of the podcasts I listen to; although"

# 6. "an average podcast last an hour" -> "an average podcast lasts an hour"
FR "an average podcast last an hour" "an average podcast lasts an hour"

# 7. " Not really" -> " It is complicated"
FR " Not really" " It is complicated"

# 8. "Some of the artists I listen to the most made the chat, but it was"
FR "Some of the artists I listen to the most made the chat, but it was" "Some of the top artists I listen to made the chart, but it was"

# 9. Big naval/solid paragraph rewrite
FR "The top song (in this case podcast), How to Get Rich, accounted for 40% of the time I spent listening to naval and solid, although not a drake song charted as top 3, goes to show how much I like drake." "The top song (in this case podcast), How to Get Rich, accounted for 40% of the time I spent listening to the artist Naval. Although Young Stoner Boy featured drake in the song Solid, it charted as top 3, goes to show how much I like drake. Lemon Pepper Freestyle by Drake was in the chart too. I never thought my most listened to song would be a podcast, but it was expected given the average listening time of a podcast."

Write-Output "text edits done"
